# Corrigindo erros de 'id' duplicado
# Fix duplicate 'id' column values: renumber existing products (rows 2-3),
# add a third sample product (row 4) that was sharing id=2/SKU with row 3,
# and append three new placeholder products (rows 5-7) each with their own
# unique id (4, 5, 6).
#
# NOTE: string values are written with a leading "'" (apostrophe) so the
# host keeps numeric-looking / boolean-looking text ("299.99", "true", ...)
# as literal text instead of auto-coercing it to a number/boolean cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (id=1) -------------------------------------------------------
# Product name + SKU text were corrected; is_available flipped to false.
$ws.Cells.Item(2, 1).Value  = 1
$ws.Cells.Item(2, 2).Value  = "'M"
$ws.Cells.Item(2, 3).Value  = "'Camisa 233323"
$ws.Cells.Item(2, 4).Value  = "'299.99"
$ws.Cells.Item(2, 5).Value  = "'399.99"
$ws.Cells.Item(2, 6).Value  = "'SKU-hdhdhdhdh"
$ws.Cells.Item(2, 7).Value  = "'Camisa-azul"
$ws.Cells.Item(2, 8).Value  = "'https://i.imgur.com/NIdnVcg.jpeg"
$ws.Cells.Item(2, 9).Value  = "'https://i.imgur.com/pKSmw8F.jpg"
$ws.Cells.Item(2, 10).Value = "'https://i.imgur.com/vmnEV31.jpg"
$ws.Cells.Item(2, 11).Value = "'https://i.imgur.com/5TN2v5a.jpg"
$ws.Cells.Item(2, 12).Value = "'https://i.imgur.com/hTBFx7g.jpg"
$ws.Cells.Item(2, 13).Value = "'false"
$ws.Cells.Item(2, 14).Value = "'false"

# --- Row 3 (id=2) -------------------------------------------------------
$ws.Cells.Item(3, 1).Value  = 2
$ws.Cells.Item(3, 2).Value  = "'M"
$ws.Cells.Item(3, 3).Value  = "'camisa aslkdla"
$ws.Cells.Item(3, 4).Value  = "'899.89"
$ws.Cells.Item(3, 5).Value  = "'99.99"
$ws.Cells.Item(3, 6).Value  = "'SKU-245345"
$ws.Cells.Item(3, 7).Value  = "'Camisa-rosa"
$ws.Cells.Item(3, 8).Value  = "'https://i.imgur.com/NIdnVcg.jpeg"
$ws.Cells.Item(3, 9).Value  = "'https://i.imgur.com/pKSmw8F.jpg"
$ws.Cells.Item(3, 10).Value = "'https://i.imgur.com/vmnEV31.jpg"
$ws.Cells.Item(3, 11).Value = "'https://i.imgur.com/5TN2v5a.jpg"
$ws.Cells.Item(3, 12).Value = "'https://i.imgur.com/hTBFx7g.jpg"
$ws.Cells.Item(3, 13).Value = "'true"
$ws.Cells.Item(3, 14).Value = "'true"

# --- Row 4 (id=3, new row) ----------------------------------------------
$ws.Cells.Item(4, 1).Value  = 3
$ws.Cells.Item(4, 2).Value  = "'F"
$ws.Cells.Item(4, 3).Value  = "'camisa emily2"
$ws.Cells.Item(4, 4).Value  = "'899.89"
$ws.Cells.Item(4, 5).Value  = "'99.99"
$ws.Cells.Item(4, 6).Value  = "'SKU-245345"
$ws.Cells.Item(4, 7).Value  = "'Camisa-preta"
$ws.Cells.Item(4, 8).Value  = "'https://i.imgur.com/NIdnVcg.jpeg"
$ws.Cells.Item(4, 9).Value  = "'https://i.imgur.com/pKSmw8F.jpg"
$ws.Cells.Item(4, 10).Value = "'https://i.imgur.com/vmnEV31.jpg"
$ws.Cells.Item(4, 11).Value = "'https://i.imgur.com/5TN2v5a.jpg"
$ws.Cells.Item(4, 12).Value = "'https://i.imgur.com/hTBFx7g.jpg"
$ws.Cells.Item(4, 13).Value = "'false"
$ws.Cells.Item(4, 14).Value = "'true"

# --- Row 5 (id=4, new row) ----------------------------------------------
$ws.Cells.Item(5, 1).Value  = 4
$ws.Cells.Item(5, 2).Value  = "'c "
$ws.Cells.Item(5, 3).Value  = "'ccassmdiasmd"
$ws.Cells.Item(5, 4).Value  = 34534
$ws.Cells.Item(5, 5).Value  = 45345
$ws.Cells.Item(5, 6).Value  = "'fsdfsdf"
$ws.Cells.Item(5, 7).Value  = "'sdfsdfsdf"
$ws.Cells.Item(5, 8).Value  = "'sdfsdfsd"
$ws.Cells.Item(5, 9).Value  = "'fsdfsdf"
$ws.Cells.Item(5, 10).Value = "'fsdfsdfsdf"
$ws.Cells.Item(5, 11).Value = "'sdfsdfsdfs"
$ws.Cells.Item(5, 12).Value = "'fsdfsdfsdf"
$ws.Cells.Item(5, 13).Value = "'true"
$ws.Cells.Item(5, 14).Value = "'false"

# --- Row 6 (id=5, new row) ----------------------------------------------
$ws.Cells.Item(6, 1).Value  = 5
$ws.Cells.Item(6, 2).Value  = "'c "
$ws.Cells.Item(6, 3).Value  = "'asdasda"
$ws.Cells.Item(6, 4).Value  = 456456
$ws.Cells.Item(6, 5).Value  = 5645
$ws.Cells.Item(6, 6).Value  = "'sfsdfsdf"
$ws.Cells.Item(6, 7).Value  = "'sdfsdfsdf"
$ws.Cells.Item(6, 8).Value  = "'sdfsdfsdf"
$ws.Cells.Item(6, 9).Value  = "'sdfsdfsdf"
$ws.Cells.Item(6, 10).Value = "'sdfsdfsdfsdfs"
$ws.Cells.Item(6, 11).Value = "'sdfsdfsdfsdf"
$ws.Cells.Item(6, 12).Value = "'sdfsdfsdfsdf"
$ws.Cells.Item(6, 13).Value = "'false"
$ws.Cells.Item(6, 14).Value = "'true"

# --- Row 7 (id=6, new row) ----------------------------------------------
$ws.Cells.Item(7, 1).Value  = 6
$ws.Cells.Item(7, 2).Value  = "'M"
$ws.Cells.Item(7, 3).Value  = "'sadasdsdfa"
$ws.Cells.Item(7, 4).Value  = 7869789
$ws.Cells.Item(7, 5).Value  = 67867
$ws.Cells.Item(7, 6).Value  = "'fghfghfgh"
$ws.Cells.Item(7, 7).Value  = "'fghfghfg"
$ws.Cells.Item(7, 8).Value  = "'fghfghfgh"
$ws.Cells.Item(7, 9).Value  = "'fghfghfghfg"
$ws.Cells.Item(7, 10).Value = "'fghfghfghfgh"
$ws.Cells.Item(7, 11).Value = "'fghfghfghfghfghfg"
$ws.Cells.Item(7, 12).Value = "'hfghfghfghfg"
$ws.Cells.Item(7, 13).Value = "'true"
$ws.Cells.Item(7, 14).Value = "'false"

# --- View bookkeeping: match the author's final selection/scroll --------
$ws.Range("N2").Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollColumn = 9
    $win.ScrollRow = 1
}
